$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date value from 45184 to 45186 for rows 2-20
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 3).Value = 45186
}

# Update HYPERLINK formulas in columns S, T, V, W, X, Y for rows 2-4 to add the
# friendly-name second argument (the beteckning in column A of that row)
$cols = @(19, 20, 22, 23, 24, 25)  # S, T, V, W, X, Y
for ($r = 2; $r -le 4; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    foreach ($c in $cols) {
        $cell = $ws.Cells.Item($r, $c)
        $f = $cell.Formula
        if ($f -match '^=HYPERLINK\("([^"]*)"\)$') {
            $url = $matches[1]
            $cell.Formula = '=HYPERLINK("' + $url + '", "' + $name + '")'
        }
    }
}
